$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.371231
$ws.Range("H2").Value = 25.113693
$ws.Range("I2").Value = 0.1018690981095697
$ws.Range("J2").Value = 0.1018690981095697
$ws.Range("M2").Value = 82.43338033333333
$ws.Range("N2").Value = 247.300141
$ws.Range("O2").Value = 0.3670006993429558
$ws.Range("P2").Value = 0.3670006993429557
$ws.Range("Q2").Value = 690.0688688811903
$ws.Range("R2").Value = 6210.619819930713
$ws.Range("S2").Value = 0.03738603024764826
$ws.Range("T2").Value = 0.03738603024764825
$ws.Range("G3").Value = 8.371231
$ws.Range("H3").Value = 25.113693
$ws.Range("I3").Value = 0.1018690981095697
$ws.Range("J3").Value = 0.1018690981095697
$ws.Range("O3").Value = 0.3956886215996139
$ws.Range("P3").Value = 0.3956886215996139
$ws.Range("Q3").Value = 744.0105700758902
$ws.Range("R3").Value = 6696.095130683012
$ws.Range("S3").Value = 0.04030844301457148
$ws.Range("T3").Value = 0.04030844301457147
$ws.Range("G4").Value = 8.371231
$ws.Range("H4").Value = 25.113693
$ws.Range("I4").Value = 0.1018690981095697
$ws.Range("J4").Value = 0.1018690981095697
$ws.Range("M4").Value = 42.93483766666667
$ws.Range("N4").Value = 128.804513
$ws.Range("O4").Value = 0.1911496942879982
$ws.Range("P4").Value = 0.1911496942879981
$ws.Range("Q4").Value = 359.4174440551677
$ws.Range("R4").Value = 3234.756996496509
$ws.Range("S4").Value = 0.01947224696103834
$ws.Range("T4").Value = 0.01947224696103834
$ws.Range("G5").Value = 8.371231
$ws.Range("H5").Value = 25.113693
$ws.Range("I5").Value = 0.1018690981095697
$ws.Range("J5").Value = 0.1018690981095697
$ws.Range("M5").Value = 10.368389
$ws.Range("N5").Value = 31.105167
$ws.Range("O5").Value = 0.04616098476943217
$ws.Range("P5").Value = 0.04616098476943217
$ws.Range("Q5").Value = 86.796179416859
$ws.Range("R5").Value = 781.165614751731
$ws.Range("S5").Value = 0.004702377886311639
$ws.Range("T5").Value = 0.004702377886311638
$ws.Range("I6").Value = 0.683327746432814
$ws.Range("J6").Value = 0.683327746432814
$ws.Range("M6").Value = 82.43338033333333
$ws.Range("N6").Value = 247.300141
$ws.Range("O6").Value = 0.3670006993429558
$ws.Range("P6").Value = 0.3670006993429557
$ws.Range("Q6").Value = 4628.913122886747
$ws.Range("R6").Value = 41660.21810598073
$ws.Range("S6").Value = 0.2507817608212887
$ws.Range("T6").Value = 0.2507817608212887
$ws.Range("I7").Value = 0.683327746432814
$ws.Range("J7").Value = 0.683327746432814
$ws.Range("O7").Value = 0.3956886215996139
$ws.Range("P7").Value = 0.3956886215996139
$ws.Range("S7").Value = 0.2703850140867707
$ws.Range("T7").Value = 0.2703850140867706
$ws.Range("I8").Value = 0.683327746432814
$ws.Range("J8").Value = 0.683327746432814
$ws.Range("M8").Value = 42.93483766666667
$ws.Range("N8").Value = 128.804513
$ws.Range("O8").Value = 0.1911496942879982
$ws.Range("P8").Value = 0.1911496942879981
$ws.Range("Q8").Value = 2410.936354915935
$ws.Range("R8").Value = 21698.42719424341
$ws.Range("S8").Value = 0.1306178898291391
$ws.Range("T8").Value = 0.1306178898291391
$ws.Range("I9").Value = 0.683327746432814
$ws.Range("J9").Value = 0.683327746432814
$ws.Range("M9").Value = 10.368389
$ws.Range("N9").Value = 31.105167
$ws.Range("O9").Value = 0.04616098476943217
$ws.Range("P9").Value = 0.04616098476943217
$ws.Range("Q9").Value = 582.2201116977277
$ws.Range("R9").Value = 5239.981005279549
$ws.Range("S9").Value = 0.03154308169561554
$ws.Range("T9").Value = 0.03154308169561553
$ws.Range("G10").Value = 16.77784
$ws.Range("H10").Value = 50.33351999999999
$ws.Range("I10").Value = 0.2041687093602677
$ws.Range("J10").Value = 0.2041687093602677
$ws.Range("M10").Value = 82.43338033333333
$ws.Range("N10").Value = 247.300141
$ws.Range("O10").Value = 0.3670006993429558
$ws.Range("P10").Value = 0.3670006993429557
$ws.Range("Q10").Value = 1383.054065891813
$ws.Range("R10").Value = 12447.48659302632
$ws.Range("S10").Value = 0.07493005911916692
$ws.Range("T10").Value = 0.07493005911916691
$ws.Range("G11").Value = 16.77784
$ws.Range("H11").Value = 50.33351999999999
$ws.Range("I11").Value = 0.2041687093602677
$ws.Range("J11").Value = 0.2041687093602677
$ws.Range("O11").Value = 0.3956886215996139
$ws.Range("P11").Value = 0.3956886215996139
$ws.Range("Q11").Value = 1491.165433499813
$ws.Range("R11").Value = 13420.48890149832
$ws.Range("S11").Value = 0.08078723518053652
$ws.Range("T11").Value = 0.0807872351805365
$ws.Range("G12").Value = 16.77784
$ws.Range("H12").Value = 50.33351999999999
$ws.Range("I12").Value = 0.2041687093602677
$ws.Range("J12").Value = 0.2041687093602677
$ws.Range("M12").Value = 42.93483766666667
$ws.Range("N12").Value = 128.804513
$ws.Range("O12").Value = 0.1911496942879982
$ws.Range("P12").Value = 0.1911496942879981
$ws.Range("Q12").Value = 720.3538367973067
$ws.Range("R12").Value = 6483.18453117576
$ws.Range("S12").Value = 0.03902678637739031
$ws.Range("T12").Value = 0.03902678637739031
$ws.Range("G13").Value = 16.77784
$ws.Range("H13").Value = 50.33351999999999
$ws.Range("I13").Value = 0.2041687093602677
$ws.Range("J13").Value = 0.2041687093602677
$ws.Range("M13").Value = 10.368389
$ws.Range("N13").Value = 31.105167
$ws.Range("O13").Value = 0.04616098476943217
$ws.Range("P13").Value = 0.04616098476943217
$ws.Range("Q13").Value = 173.95917169976
$ws.Range("R13").Value = 1565.63254529784
$ws.Range("S13").Value = 0.00942462868317394
$ws.Range("T13").Value = 0.009424628683173938
$ws.Range("G14").Value = 0.8738999999999999
$ws.Range("H14").Value = 2.6217
$ws.Range("I14").Value = 0.01063444609734852
$ws.Range("J14").Value = 0.01063444609734852
$ws.Range("M14").Value = 82.43338033333333
$ws.Range("N14").Value = 247.300141
$ws.Range("O14").Value = 0.3670006993429558
$ws.Range("P14").Value = 0.3670006993429557
$ws.Range("Q14").Value = 72.03853107329999
$ws.Range("R14").Value = 648.3467796596999
$ws.Range("S14").Value = 0.003902849154851874
$ws.Range("T14").Value = 0.003902849154851873
$ws.Range("G15").Value = 0.8738999999999999
$ws.Range("H15").Value = 2.6217
$ws.Range("I15").Value = 0.01063444609734852
$ws.Range("J15").Value = 0.01063444609734852
$ws.Range("O15").Value = 0.3956886215996139
$ws.Range("P15").Value = 0.3956886215996139
$ws.Range("Q15").Value = 77.66968050329999
$ws.Range("R15").Value = 699.0271245296999
$ws.Range("S15").Value = 0.00420792931773523
$ws.Range("T15").Value = 0.00420792931773523
$ws.Range("G16").Value = 0.8738999999999999
$ws.Range("H16").Value = 2.6217
$ws.Range("I16").Value = 0.01063444609734852
$ws.Range("J16").Value = 0.01063444609734852
$ws.Range("M16").Value = 42.93483766666667
$ws.Range("N16").Value = 128.804513
$ws.Range("O16").Value = 0.1911496942879982
$ws.Range("P16").Value = 0.1911496942879981
$ws.Range("Q16").Value = 37.5207546369
$ws.Range("R16").Value = 337.6867917321
$ws.Range("S16").Value = 0.002032771120430365
$ws.Range("T16").Value = 0.002032771120430365
$ws.Range("G17").Value = 0.8738999999999999
$ws.Range("H17").Value = 2.6217
$ws.Range("I17").Value = 0.01063444609734852
$ws.Range("J17").Value = 0.01063444609734852
$ws.Range("M17").Value = 10.368389
$ws.Range("N17").Value = 31.105167
$ws.Range("O17").Value = 0.04616098476943217
$ws.Range("P17").Value = 0.04616098476943217
$ws.Range("Q17").Value = 9.060935147099999
$ws.Range("R17").Value = 81.5484163239
$ws.Range("S17").Value = 0.0004908965043310525
$ws.Range("T17").Value = 0.0004908965043310525
